# tempocctraits.xlsx - "temp occ traits saved for ms resubmission"
#
# Adds a new worksheet "Sheet1" (placed after "combined") that holds a
# side-by-side summary pasted from three query-table refreshes:
#   B:C -> "Original results; logit transform"          (Estimate / p-value)
#   D:E -> "Including median occupancy; logit transform" (Estimate / p-value)
#   F:G -> "Including median occupancy; arcsin transform"(Estimate / p-value)
# Also adds the three sheet-scoped named ranges for Sheet1, switches the
# active tab to Sheet1, and switches a few of the smallest p-values on the
# "combined" sheet to scientific notation.

$wb = $excel.ActiveWorkbook
$combined = $wb.Worksheets.Item("combined")

# ---------------------------------------------------------------------
# 1. A few p-values on "combined" get reformatted to scientific notation
# ---------------------------------------------------------------------
$combined.Range("E5").NumberFormat = "0.000E+00"
$combined.Range("E7").NumberFormat = "0.000E+00"
$combined.Range("E15").NumberFormat = "0.000E+00"

# ---------------------------------------------------------------------
# 2. Add the new sheet, right after "combined"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $combined)
$ws.Name = "Sheet1"

# Column widths
$ws.Columns.Item(1).ColumnWidth = 20.140625
$ws.Columns.Item(2).ColumnWidth = 8.85546875
$ws.Columns.Item(3).ColumnWidth = 9.140625
$ws.Columns.Item(4).ColumnWidth = 8.5703125
$ws.Columns.Item(5).ColumnWidth = 9
$ws.Columns.Item(6).ColumnWidth = 9.28515625
$ws.Columns.Item(7).ColumnWidth = 9.7109375

# Row 1 height (group headers)
$ws.Rows.Item(1).RowHeight = 45.75

# ---------------------------------------------------------------------
# 3. Group header row (row 1) + merges
# ---------------------------------------------------------------------
$ws.Range("B1").Value = "Original results; logit transform"
$ws.Range("D1").Value = "Including median occupancy; logit transform"
$ws.Range("F1").Value = "Including median occupancy; arcsin transform"

$ws.Range("B1:C1").Merge()
$ws.Range("D1:E1").Merge()
$ws.Range("F1:G1").Merge()

# ---------------------------------------------------------------------
# 4. Column header row (row 2): Estimate / p-value
# ---------------------------------------------------------------------
$ws.Range("B2").Value = "Estimate"
$ws.Range("C2").Value = "p-value"
$ws.Range("D2").Value = "Estimate"
$ws.Range("E2").Value = "p-value"
$ws.Range("F2").Value = "Estimate"
$ws.Range("G2").Value = "p-value"

# ---------------------------------------------------------------------
# 5. Row labels (column A), copied from "combined"
# ---------------------------------------------------------------------
$labels = @(
  "Intercept (Granivore/ Neotropical Migrant)",
  "Log10(Focal Range Size)",
  "Median Temporal Occupancy",
  "Competitor Range Overlap",
  "Temperature",
  "Precipitation",
  "Elevation",
  "NDVI",
  "Insectivore/Omnivore",
  "Insectivore",
  "Omnivore",
  "Resident",
  "Short-distance Migrant"
)
for ($i = 0; $i -lt $labels.Length; $i++) {
  $ws.Cells.Item(3 + $i, 1).Value = $labels[$i]
}

# ---------------------------------------------------------------------
# 6. Data values
# ---------------------------------------------------------------------
$little = "< 2e-16"

# B:C  Original results; logit transform  (Estimate, p-value)
$ws.Range("B3").Value = 15.9996098
$ws.Range("C3").Value = 0.000000000000219
$ws.Range("B4").Value = -2.3296407
$ws.Range("C4").Value = $little
$ws.Range("B5").Value = $null
$ws.Range("C5").Value = $null
$ws.Range("B6").Value = 0.049043499999999997
$ws.Range("C6").Value = 0.00076400000000000003
$ws.Range("B7").Value = -0.0093626999999999998
$ws.Range("C7").Value = 0.020129000000000001
$ws.Range("B8").Value = 0.0106243
$ws.Range("C8").Value = 0.120033
$ws.Range("B9").Value = 0.0003589
$ws.Range("C9").Value = 0.379225
$ws.Range("B10").Value = -2.0131505000000001
$ws.Range("C10").Value = 0.269507
$ws.Range("B11").Value = 1.0039781999999999
$ws.Range("C11").Value = 0.00060400000000000004
$ws.Range("B12").Value = 0.65705259999999999
$ws.Range("C12").Value = 0.024990999999999999
$ws.Range("B13").Value = -0.38520569999999998
$ws.Range("C13").Value = 0.32131599999999999
$ws.Range("B14").Value = -0.7328538
$ws.Range("C14").Value = 0.0000776
$ws.Range("B15").Value = 0.60133550000000002
$ws.Range("C15").Value = 0.0000025299999999999999

# D:E  Including median occupancy; logit transform (Estimate, p-value)
$ws.Range("D3").Value = 19.509011000000001
$ws.Range("E3").Value = $little
$ws.Range("D4").Value = -2.880493
$ws.Range("E4").Value = $little
$ws.Range("D5").Value = 2.3747796999999999
$ws.Range("E5").Value = 0.0000000000000017400000000000001
$ws.Range("D6").Value = 0.041703499999999998
$ws.Range("E6").Value = 0.0033240000000000001
$ws.Range("D7").Value = -0.0151434
$ws.Range("E7").Value = 0.00014999999999999999
$ws.Range("D8").Value = 0.023181799999999999
$ws.Range("E8").Value = 0.00070699999999999995
$ws.Range("D9").Value = 0.00046880000000000002
$ws.Range("E9").Value = 0.23828299999999999
$ws.Range("D10").Value = -4.5671888000000003
$ws.Range("E10").Value = 0.01137
$ws.Range("D11").Value = 0.67589580000000005
$ws.Range("E11").Value = 0.018678
$ws.Range("D12").Value = 0.51955110000000004
$ws.Range("E12").Value = 0.068964999999999999
$ws.Range("D13").Value = -0.60181609999999996
$ws.Range("E13").Value = 0.11239
$ws.Range("D14").Value = -0.45834150000000001
$ws.Range("E14").Value = 0.012406
$ws.Range("D15").Value = 0.47688150000000001
$ws.Range("E15").Value = 0.00013799999999999999

# F:G  Including median occupancy; arcsin transform (Estimate, p-value)
$ws.Range("F3").Value = 4.8540000000000001
$ws.Range("G3").Value = $little
$ws.Range("F4").Value = -0.57030000000000003
$ws.Range("G4").Value = $little
$ws.Range("F5").Value = 0.41689999999999999
$ws.Range("G5").Value = 0.00000000000059599999999999998
$ws.Range("F6").Value = 0.0080440000000000008
$ws.Range("G6").Value = 0.003591
$ws.Range("F7").Value = -0.0034759999999999999
$ws.Range("G7").Value = 0.0000078399999999999995
$ws.Range("F8").Value = 0.0045640000000000003
$ws.Range("G8").Value = 0.00060800000000000003
$ws.Range("F9").Value = 0.000025680000000000001
$ws.Range("G9").Value = 0.73970899999999995
$ws.Range("F10").Value = -1.0640000000000001
$ws.Range("G10").Value = 0.0024390000000000002
$ws.Range("F11").Value = 0.18490000000000001
$ws.Range("G11").Value = 0.00095000000000000001
$ws.Range("F12").Value = 0.15870000000000001
$ws.Range("G12").Value = 0.0043299999999999996
$ws.Range("F13").Value = -0.11169999999999999
$ws.Range("G13").Value = 0.12967300000000001
$ws.Range("F14").Value = -0.052800000000000003
$ws.Range("G14").Value = 0.138269
$ws.Range("F15").Value = 0.1065
$ws.Range("G15").Value = 0.0000124

# ---------------------------------------------------------------------
# 7. Formatting
# ---------------------------------------------------------------------

# Whole used range: Segoe UI 10 pt + full thin box border
$all = $ws.Range("A1:G15")
$all.Font.Name = "Segoe UI"
$all.Font.Size = 10
$all.Borders.LineStyle = 1
$all.Borders.Weight = 2

# Body number format (row labels + data): 0.000
$ws.Range("B3:G15").NumberFormat = "0.000"

# Column A wraps text
$ws.Range("A2:A15").WrapText = $true

# Re-apply scientific notation to the small p-values
$ws.Range("C3").NumberFormat = "0.000E+00"
$ws.Range("E5").NumberFormat = "0.000E+00"
$ws.Range("G5").NumberFormat = "0.000E+00"
$ws.Range("E7").NumberFormat = "0.000E+00"
$ws.Range("G7").NumberFormat = "0.000E+00"
$ws.Range("B9").NumberFormat = "0.000E+00"
$ws.Range("C14").NumberFormat = "0.000E+00"
$ws.Range("C15").NumberFormat = "0.000E+00"
$ws.Range("E15").NumberFormat = "0.000E+00"
$ws.Range("G15").NumberFormat = "0.000E+00"

# Group header row (row 1): centered, wrapped
$ws.Range("A1:G1").WrapText = $true
$ws.Range("B1:G1").HorizontalAlignment = -4108
$ws.Range("F1:G1").VerticalAlignment = -4108

# ---------------------------------------------------------------------
# 8. Named ranges, scoped to the new sheet
# ---------------------------------------------------------------------
$ws.Names.Add("trait_mod_output_arcsine", "=Sheet1!`$F`$2:`$G`$15")
$ws.Names.Add("trait_mod_output_logit_1", "=Sheet1!`$D`$2:`$E`$15")
$ws.Names.Add("trait_mod_output_og_logit", "=Sheet1!`$A`$2:`$C`$15")

# ---------------------------------------------------------------------
# 9. Selection / active sheet (Sheet1 becomes the active tab)
# ---------------------------------------------------------------------
$combined.Range("A1:O15").Select()
$ws.Range("M11").Select()
$ws.Activate()
